# Regenerate save_data column "K" (column G) with new strike-out-derived
# values (K replaces the old "Strike#" derived figures), per recalculated
# std/mean and s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 2
    16 = 2
    17 = 2
    18 = 3
    19 = 1
    20 = 2
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
